# "All analysis page added and also download functionality of
#  'Patient Admitted in Watchlist Hospitals'"
#
# The underlying data refresh collapses every per-row "Number of Surgeons"
# (col D) / "Number of OT" (col E) count down to a flat 1 (rows 2-17) and
# the previously selected/scrolled view is updated to highlight the new
# "Number of OT" column (E2:E17) that now feeds the watchlist download.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Number of Surgeons) and E (Number of OT) -> 1 for every
# data row (2 through 17).
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 4).Value = 1
    $ws.Cells.Item($row, 5).Value = 1
}

# Rows 16 & 17 lose their explicit custom row height (back to sheet default).
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(17).AutoFit()

# Update the view: scroll so row 4 is the top visible row, and select the
# full "Number of OT" column range (E2:E17) with E2 as the active cell.
$null = $ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("E2:E17").Select()
